# Razer Key Map Visualized - update keymapping
# - Removed "Open Spotify" text from button 4 (D2), replaced with "Hold Spacebar Down"
# - Added macro to rapidly click M1 in button 3 (C2), which also picks up the
#   green highlight + wrapped text formatting used by the other macro buttons
# - Row 3 number labels above C2/D2 (C3/D3) get wrap-text turned on to match

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Button 3 (C2): was a plain number "3" (bold font + thin border already applied),
# now becomes macro text and picks up the same green fill + wrapped text used by
# the other macro buttons (e.g. D2's "Open Spotify" -> green cell)
$ws.Range("C2").Value = "3`nRapidly Click M1"
$ws.Range("C2").WrapText = $true
$ws.Range("C2").Interior.Color = 5287936

# Button 4 (D2): swap out the Spotify macro for the new "hold spacebar" macro,
# keep its existing green/wrapped style as-is
$ws.Range("D2").Value = "4`nHold Spacebar Down"

# Turn on wrap text for the plain numbers sitting above the new macro text (C3/D3)
$ws.Range("C3").WrapText = $true
$ws.Range("D3").WrapText = $true

# Reflect the last clicked cell in the saved view (cosmetic, matches authoring session)
$ws.Range("G11").Select()
